$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.43869999999999
$ws.Range("B3").Value = 6.194900000000003
$ws.Range("B14").Value = 5.911500000000001
$ws.Range("B21").Value = 9.703700000000003
$ws.Range("B23").Value = 9.218100000000002
$ws.Range("B25").Value = 5.336300000000002
$ws.Range("E25").Value = 17.2288
$ws.Range("B26").Value = 4.313500000000004
$ws.Range("E27").Value = 16.78109999999998
$ws.Range("B29").Value = 5.047100000000003
$ws.Range("E31").Value = 15.9783
$ws.Range("E39").Value = 16.22319999999999
$ws.Range("E48").Value = 17.4088
$ws.Range("E51").Value = 17.18470000000001
$ws.Range("E52").Value = 17.06930000000001
$ws.Range("B53").Value = 5.200100000000003
$ws.Range("E55").Value = 16.632
$ws.Range("E56").Value = 16.08050000000001
$ws.Range("B57").Value = 4.849099999999995
$ws.Range("E57").Value = 16.7187
$ws.Range("B59").Value = 5.2326
$ws.Range("B69").Value = 5.2296
$ws.Range("E73").Value = 17.17150000000001
$ws.Range("B79").Value = 9.457400000000007
$ws.Range("B83").Value = 5.173599999999998
$ws.Range("E89").Value = 17.25190000000002
$ws.Range("E90").Value = 16.83249999999999
$ws.Range("B91").Value = 5.194900000000001
$ws.Range("E92").Value = 18.56730000000001
$ws.Range("B93").Value = 5.763699999999999
